# Weekly price-sheet update: insert the newest "Achicoria" market reading
# as a new row 24, pushing the previously-existing rows 24-33 down to 25-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24..33 down to 25..34 by inserting a new blank row at 24.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly record.
$ws.Cells.Item(24, 1).Value = 9
$ws.Cells.Item(24, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44943
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = 100112010
$ws.Cells.Item(24, 7).Value = "Achicoria"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 70
$ws.Cells.Item(24, 11).Value = 7000
$ws.Cells.Item(24, 12).Value = 7000
$ws.Cells.Item(24, 13).Value = 7000
$ws.Cells.Item(24, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(24, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(24, 16).Value = 438
$ws.Cells.Item(24, 17).Value = 16
$ws.Cells.Item(24, 18).Value = "Hortaliza"
